# Weekly update: insert a new week's record at row 18 (shifts existing
# data rows 18-61 down to 19-62) for the Albahaca / Terminal La Palmera
# de La Serena sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 18; rows 18-61 shift to 19-62.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with this week's values.
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = "Terminal La Palmera de La Serena"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44525
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 100112052
$ws.Range("G18").Value = "Albahaca"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 720
$ws.Range("K18").Value = 3000
$ws.Range("L18").Value = 4000
$ws.Range("M18").Value = 3500
$ws.Range("N18").Value = "$/paquete"
$ws.Range("O18").Value = "Región de Arica y Parinacota"
$ws.Range("P18").Value = 3500
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
